# Weekly data refresh: insert a new price-record row for
# "Feria Lagunitas de Puerto Montt - Zapallo italiano" ahead of the
# existing history, shifting all later rows down by one.
#
# Before the edit, row 201 holds the record for fecha=44390 (2021-07-13).
# A brand-new record for fecha=44736 (2022-06-24) is inserted above it,
# so everything that used to be row 201..250 becomes row 202..251, and the
# used range grows from A1:R250 to A1:R251.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 201 (pushes old 201..250 down to 202..251).
$ws.Rows.Item(201).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(201, 1).Value = 4
$ws.Cells.Item(201, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(201, 3).Value = "Los Lagos"
$ws.Cells.Item(201, 4).Value = 44736
$ws.Cells.Item(201, 5).Value = 10
$ws.Cells.Item(201, 6).Value = 100112032
$ws.Cells.Item(201, 7).Value = "Zapallo italiano"
$ws.Cells.Item(201, 8).Value = "Sin especificar"
$ws.Cells.Item(201, 9).Value = "Primera"
$ws.Cells.Item(201, 10).Value = 200
$ws.Cells.Item(201, 11).Value = 20000
$ws.Cells.Item(201, 12).Value = 21000
$ws.Cells.Item(201, 13).Value = 20500
$ws.Cells.Item(201, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(201, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(201, 16).Value = 410
$ws.Cells.Item(201, 17).Value = 50
$ws.Cells.Item(201, 18).Value = "Hortaliza"
